$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Insert a new row at the top (old row 1 header becomes row 2, etc.)
$ws.Rows.Item(1).Insert()

# New title row: D1 = Preset1, G1 = Preset2
$ws.Range("D1").Value = "Preset1"
$ws.Range("G1").Value = "Preset2"

# New header cell for the GrainRate column
$ws.Range("G2").Value = "GrainRate"

# Fix the bug: MorphX row's Freq (D) parameter count was 0, should be 10
$ws.Range("D3").Value = 10

# New GrainRate column values (all 0) for the 4 data rows
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 0
